# "crear archivo crud de los productos"
# Add a first data row to the "productos" sheet and make it the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("productos")

# Switch focus to the productos sheet (this becomes the workbook's active tab).
$ws.Activate()

# New product row: id, nombre, categoria, precio, cantidad
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "roa"
$ws.Range("C2").Value = "cereal"
$ws.Range("D2").Value = 2500
$ws.Range("E2").Value = 4

# Leave the selection on A3, like in the saved workbook.
$ws.Range("A3").Select()
